$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 120
$ws1.Range("F4").Value = 8513
$ws1.Range("F5").Value = 568
$ws1.Range("F6").Value = 7538
$ws1.Range("F10").Value = 771
$ws1.Range("F11").Value = 127
$ws1.Range("F12").Value = 207
$ws1.Range("F13").Value = 12612
$ws1.Range("F15").Value = 111
$ws1.Range("F16").Value = 2743
$ws1.Range("F18").Value = 5944
$ws1.Range("F20").Value = 3147
$ws1.Range("F22").Value = 143
$ws1.Range("F26").Value = 46
$ws1.Range("F27").Value = 113
$ws1.Range("F28").Value = 3464
$ws1.Range("F30").Value = 2573
$ws1.Range("F32").Value = 1802
$ws1.Range("F34").Value = 169
$ws1.Range("F35").Value = 6271
$ws1.Range("F38").Value = 1273
$ws1.Range("F40").Value = 974
$ws1.Range("F41").Value = 187
$ws1.Range("F44").Value = 1126
$ws1.Range("F45").Value = 103
$ws1.Range("F46").Value = 1140
$ws1.Range("F47").Value = 1658
$ws1.Range("F48").Value = 40
$ws1.Range("F49").Value = 129
$ws1.Range("F50").Value = 1149

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 217
$ws2.Range("F18").Value = 921
$ws2.Range("F20").Value = 84

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 380
$ws3.Range("F3").Value = 541

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 120
$ws4.Range("F5").Value = 380
$ws4.Range("F6").Value = 8513
$ws4.Range("F7").Value = 568
$ws4.Range("F8").Value = 7538
$ws4.Range("F12").Value = 771
$ws4.Range("F14").Value = 207
$ws4.Range("F15").Value = 217
$ws4.Range("F16").Value = 12612
$ws4.Range("F18").Value = 2743
$ws4.Range("F19").Value = 5944
$ws4.Range("F20").Value = 3147
$ws4.Range("F22").Value = 143
$ws4.Range("F25").Value = 46
$ws4.Range("F26").Value = 113
$ws4.Range("F27").Value = 3464
$ws4.Range("F29").Value = 2573
$ws4.Range("F31").Value = 1802
$ws4.Range("F32").Value = 169
$ws4.Range("F33").Value = 6271
$ws4.Range("F34").Value = 84
$ws4.Range("F38").Value = 1273
$ws4.Range("F40").Value = 974
$ws4.Range("F41").Value = 187
$ws4.Range("F44").Value = 1126
$ws4.Range("F45").Value = 103
$ws4.Range("F46").Value = 1140
$ws4.Range("F47").Value = 1658
$ws4.Range("F48").Value = 40
$ws4.Range("F49").Value = 129
$ws4.Range("F50").Value = 1149
